# Generate Report for Handback
# Update the "handback" timestamps for the f1934dc0-... entry across the
# Overview, zh-cn and de-de sheets, reflecting that the report has now
# been (re)generated.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Sheets.Item("Overview")
$ws1.Range("G3").Value = "2016-08-26 16:48:58"

$ws2 = $wb.Sheets.Item("zh-cn")
$ws2.Range("H3").Value = "2016-08-26 16:48:54"
$ws2.Range("K3").Value = "2016-08-26 16:49:15"

$ws3 = $wb.Sheets.Item("de-de")
$ws3.Range("H3").Value = "2016-08-26 16:48:58"
$ws3.Range("K3").Value = "2016-08-26 16:49:22"
